# Update cryptocurrency price list (prices + 1h volume change %)
# for data that refreshed on Sun Mar 19 16:43:53 UTC 2023.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the whole Price column as Text so numeric-looking strings
# (e.g. "337.10", "0.9981") are preserved verbatim instead of being
# auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Cells.Item(2, 4).Value = "27.751.99"
$ws.Cells.Item(2, 5).Value = "  +1.83%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "1.808.63"
$ws.Cells.Item(3, 5).Value = "  +1.39%  "

# Row 4
$ws.Cells.Item(4, 4).Value = "0.9981"
$ws.Cells.Item(4, 5).Value = "  -1.17%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "337.10"
$ws.Cells.Item(5, 5).Value = "  +0.31%  "

# Row 6
$ws.Cells.Item(6, 4).Value = "0.9929"
$ws.Cells.Item(6, 5).Value = "  -1.64%  "

# Row 7
$ws.Cells.Item(7, 4).Value = "0.3931"
$ws.Cells.Item(7, 5).Value = "  +4.10%  "

# Row 8
$ws.Cells.Item(8, 4).Value = "0.3482"
$ws.Cells.Item(8, 5).Value = "  +1.70%  "

# Row 9
$ws.Cells.Item(9, 4).Value = "48.44"
$ws.Cells.Item(9, 5).Value = "  -1.18%  "

# Row 10
$ws.Cells.Item(10, 4).Value = "1.207"
$ws.Cells.Item(10, 5).Value = "  +0.77%  "

# Row 11
$ws.Cells.Item(11, 4).Value = "0.07555"
$ws.Cells.Item(11, 5).Value = "  +1.37%  "

# Row 12
$ws.Cells.Item(12, 4).Value = "0.9913"
$ws.Cells.Item(12, 5).Value = "  -1.63%  "

# Row 13
$ws.Cells.Item(13, 4).Value = "22.27"
$ws.Cells.Item(13, 5).Value = "  +2.04%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "6.537"
$ws.Cells.Item(14, 5).Value = "  +1.10%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "1.805.49"
$ws.Cells.Item(15, 5).Value = "  +1.28%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "7.174"
$ws.Cells.Item(16, 5).Value = "  +2.26%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "0.00001107"
$ws.Cells.Item(17, 5).Value = "  +1.55%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "0.06682"
$ws.Cells.Item(18, 5).Value = "  +0.27%  "

# Row 19
$ws.Cells.Item(19, 4).Value = "85.17"
$ws.Cells.Item(19, 5).Value = "  +1.18%  "

# Row 20
$ws.Cells.Item(20, 4).Value = "0.9979"
$ws.Cells.Item(20, 5).Value = "  -1.06%  "

# Row 21
$ws.Cells.Item(21, 4).Value = "17.80"
$ws.Cells.Item(21, 5).Value = "  +2.87%  "

# Row 22
$ws.Cells.Item(22, 4).Value = "6.586"
$ws.Cells.Item(22, 5).Value = "  +2.25%  "

# Row 23
$ws.Cells.Item(23, 4).Value = "27.779.09"
$ws.Cells.Item(23, 5).Value = "  +2.07%  "

# Row 24
$ws.Cells.Item(24, 4).Value = "12.94"
$ws.Cells.Item(24, 5).Value = "  +4.40%  "

# Row 25
$ws.Cells.Item(25, 4).Value = "2.397"
$ws.Cells.Item(25, 5).Value = "  -2.79%  "

# Row 26
$ws.Cells.Item(26, 2).Value = "LidoDAOToken"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Cells.Item(26, 4).Value = "2.554"
$ws.Cells.Item(26, 5).Value = "  +0.88%  "

# Row 27
$ws.Cells.Item(27, 2).Value = "ImmutableX"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(27, 4).Value = "1.490"
$ws.Cells.Item(27, 5).Value = "  -0.32%  "

# Row 28
$ws.Cells.Item(28, 4).Value = "21.44"
$ws.Cells.Item(28, 5).Value = "  +0.27%  "

# Row 29
$ws.Cells.Item(29, 4).Value = "155.32"
$ws.Cells.Item(29, 5).Value = "  +3.41%  "

# Row 30
$ws.Cells.Item(30, 4).Value = "2.006.31"
$ws.Cells.Item(30, 5).Value = "  +1.00%  "

# Row 31
$ws.Cells.Item(31, 4).Value = "135.83"
$ws.Cells.Item(31, 5).Value = "  +2.08%  "

# Row 32
$ws.Cells.Item(32, 4).Value = "4.030"
$ws.Cells.Item(32, 5).Value = "  -1.59%  "

# Row 33
$ws.Cells.Item(33, 4).Value = "6.110"
$ws.Cells.Item(33, 5).Value = "  +1.54%  "

# Row 34
$ws.Cells.Item(34, 4).Value = "0.08807"
$ws.Cells.Item(34, 5).Value = "  +2.10%  "

# Row 35
$ws.Cells.Item(35, 4).Value = "13.26"
$ws.Cells.Item(35, 5).Value = "  +2.18%  "

# Row 36
$ws.Cells.Item(36, 4).Value = "5.521"
$ws.Cells.Item(36, 5).Value = "  +2.07%  "

# Row 37
$ws.Cells.Item(37, 4).Value = "0.02426"
$ws.Cells.Item(37, 5).Value = "  +3.92%  "

# Row 38
$ws.Cells.Item(38, 2).Value = "TheSandbox"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Cells.Item(38, 4).Value = "0.6895"
$ws.Cells.Item(38, 5).Value = "  +0.01%  "

# Row 39
$ws.Cells.Item(39, 4).Value = "0.06521"
$ws.Cells.Item(39, 5).Value = "  +3.07%  "

# Row 40
$ws.Cells.Item(40, 2).Value = "WEMIXTOKEN"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Cells.Item(40, 4).Value = "1.611"
$ws.Cells.Item(40, 5).Value = "  -2.82%  "

# Row 41
$ws.Cells.Item(41, 4).Value = "0.2222"
$ws.Cells.Item(41, 5).Value = "  +1.23%  "

# Row 42
$ws.Cells.Item(42, 4).Value = "1.268"
$ws.Cells.Item(42, 5).Value = "  +0.33%  "

# Row 43
$ws.Cells.Item(43, 4).Value = "8.518"
$ws.Cells.Item(43, 5).Value = "  -2.94%  "

# Row 44
$ws.Cells.Item(44, 4).Value = "14.58"
$ws.Cells.Item(44, 5).Value = "  +1.94%  "

# Row 45
$ws.Cells.Item(45, 4).Value = "0.6512"
$ws.Cells.Item(45, 5).Value = "  +1.53%  "

# Row 46
$ws.Cells.Item(46, 4).Value = "0.9947"
$ws.Cells.Item(46, 5).Value = "  -1.38%  "

# Row 47
$ws.Cells.Item(47, 4).Value = "3.858"
$ws.Cells.Item(47, 5).Value = "  +0.18%  "

# Row 48
$ws.Cells.Item(48, 4).Value = "2.164"
$ws.Cells.Item(48, 5).Value = "  +2.29%  "

# Row 49
$ws.Cells.Item(49, 4).Value = "132.50"
$ws.Cells.Item(49, 5).Value = "  +2.50%  "

# Row 50
$ws.Cells.Item(50, 4).Value = "0.07228"
$ws.Cells.Item(50, 5).Value = "  +0.75%  "

# Row 51
$ws.Cells.Item(51, 4).Value = "80.58"
$ws.Cells.Item(51, 5).Value = "  +2.01%  "
